$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new shared strings must be created in this exact order
# so that the rebuilt shared-strings table lands on the indices used by the
# target workbook (Q1,R1,P1,O1,S1 - matches densee_com_v3(70), densee_com_v5(71),
# vgg16_com_3_1(72), vgg16_com_3(73 - replaces old vgg16_com_2), densee_com_v5_1(74))
$ws.Range("Q1").Value = "densee_com_v3"
$ws.Range("R1").Value = "densee_com_v5"
$ws.Range("P1").Value = "vgg16_com_3_1"
$ws.Range("O1").Value = "vgg16_com_3"
$ws.Range("S1").Value = "densee_com_v5_1"

# --- Network row (row 2)
$ws.Range("P2").Value = "vgg16"
$ws.Range("Q2").Value = "densenet161"
$ws.Range("R2").Value = "densenet161"
$ws.Range("S2").Value = "densenet161"

# --- Dataset row (row 3)
$ws.Range("P3").Value = "CIFAR10"
$ws.Range("Q3").Value = "CIFAR10"
$ws.Range("R3").Value = "CIFAR10"
$ws.Range("S3").Value = "CIFAR10"

# --- Learning Rate row (row 4)
$ws.Range("P4").Value = 0.003
$ws.Range("Q4").Value = 0.004
$ws.Range("R4").Value = 0.004
$ws.Range("S4").Value = 0.004

# --- Kernel size row (row 5)
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = "on 3-3"
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 1

# --- Sigma row (row 6)
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 2

# --- Pretrained row (row 7)
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0

# --- Column widths for the new/changed columns (P,Q,R,S => 16,17,18,19)
$ws.Columns.Item(16).ColumnWidth = 20.76
$ws.Columns.Item(17).ColumnWidth = 16.6
$ws.Columns.Item(18).ColumnWidth = 15.42
$ws.Columns.Item(19).ColumnWidth = 17.26

# --- View: scroll so column B is the left-most visible column, select Q8
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("Q8").Select()

# --- Page setup (print settings)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
